$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New crosswalk rows appended after the existing data (rows 1-23), adding
# DOP (Limnoithona sinensis / tetraspina / Oithona davisae / Oithona similis)
# juvenile + adult crosswalk entries.
$newRows = @(
    @{ Row = 24; SizeClass = "Meso"; Taxname = "Limnoithona sinensis";   Lifestage = "Juvenile" },
    @{ Row = 25; SizeClass = "Meso"; Taxname = "Limnoithona tetraspina"; Lifestage = "Juvenile" },
    @{ Row = 26; SizeClass = "Meso"; Taxname = "Oithona davisae";        Lifestage = "Juvenile" },
    @{ Row = 27; SizeClass = "Meso"; Taxname = "Oithona similis";        Lifestage = "Juvenile" },
    @{ Row = 28; SizeClass = "Meso"; Taxname = "Oithona similis";        Lifestage = "Adult" }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row

    # Copy the formatting of the last existing data row (row 23, columns A:B)
    # down onto the new row so the style matches the rest of the table.
    $ws.Range("A23:B23").Copy() | Out-Null
    $ws.Range("A" + $rowNum + ":B" + $rowNum).PasteSpecial(-4122) | Out-Null

    $ws.Range("A" + $rowNum).Value = $r.SizeClass
    $ws.Range("B" + $rowNum).Value = $r.Taxname
    $ws.Range("C" + $rowNum).Value = $r.Lifestage
}

$excel.CutCopyMode = 0

# Update the active selection to reflect where the author ended up after
# entering the new data.
$ws.Range("B29").Select() | Out-Null
